{"js": "// Updated titles and descriptions of all interviews.\n// Split the second paragraph's single \"Interviewer: ...; Interviewee: ...;\n// Transcriber: ...;  Interview takes place ...\" run into separate runs\n// joined by literal \"<br />\" markers, and wrap the first \"<br />\" in a new\n// bookmark named \"__DdeLink__256_169685483\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst infoParagraph = paragraphs.items[1];\n\nconst newText =\n  \"Interviewer: Chris Judge<br />\" +\n  \"Interviewee: Chief James Caulder<br />\" +\n  \"Transcriber: Jillian Weber<br />\" +\n  \"Interview takes place at the Pee Dee Powwow. Drum group and loudspeaker in the background.\";\n\ninfoParagraph.getRange().insertText(newText, \"Replace\");\nawait context.sync();\n\n// Locate every literal \"<br />\" marker in the (now rebuilt) paragraph.\nconst paraRange = infoParagraph.getRange();\nconst breakMarks = paraRange.search(\"<br />\", { matchCase: true });\nbreakMarks.load(\"items\");\nawait context.sync();\n\n// Dropping a (temporary) bookmark on the start and end boundary of each\n// \"<br />\" forces Word to keep the surrounding text segments as separate\n// <w:r> runs instead of re-merging them into a single run back into one\n// run of identical formatting when the package is serialised.\nconst tempBookmarkNames = [];\nfor (let i = 0; i < breakMarks.items.length; i++) {\n  const mark = breakMarks.items[i];\n\n  const startName = \"__TempRunSplit_\" + i + \"_s\";\n  mark.getRange(\"Start\").insertBookmark(startName);\n  await context.sync();\n  tempBookmarkNames.push(startName);\n\n  const endName = \"__TempRunSplit_\" + i + \"_e\";\n  mark.getRange(\"End\").insertBookmark(endName);\n  await context.sync();\n  tempBookmarkNames.push(endName);\n}\n\nfor (const name of tempBookmarkNames) {\n  context.document.deleteBookmark(name);\n}\nawait context.sync();\n\n// Re-run the search (anchors shift as runs split/bookmarks come and go)\n// and wrap the first \"<br />\" occurrence in the real, permanent bookmark.\nconst paraRange2 = infoParagraph.getRange();\nconst firstBreakMark = paraRange2.search(\"<br />\", { matchCase: true });\nfirstBreakMark.load(\"items\");\nawait context.sync();\n\nfirstBreakMark.items[0].insertBookmark(\"__DdeLink__256_169685483\");\nawait context.sync();\n", "ps1": "# Updated titles and descriptions of all interviews.\n# Split the second paragraph's single \"Interviewer: ...; Interviewee: ...;\n# Transcriber: ...;  Interview takes place ...\" run into separate runs\n# joined by literal \"<br />\" markers, and wrap the first \"<br />\" in a new\n# bookmark named \"__DdeLink__256_169685483\".\n\n$d = $word.ActiveDocument\n\n$p2 = $d.Paragraphs(2).Range\n$start = $p2.Start\n\n$newText = \"Interviewer: Chris Judge<br />Interviewee: Chief James Caulder<br />Transcriber: Jillian Weber<br />Interview takes place at the Pee Dee Powwow. Drum group and loudspeaker in the background.\"\n$p2.Text = $newText\n\n# Character offsets (relative to $start) of every run boundary in the\n# rebuilt paragraph text. Inserting (and immediately deleting) a zero-width\n# bookmark at each boundary forces Word to keep the surrounding text in\n# separate <w:r> runs instead of re-merging them into one run on save.\n$boundaries = @(24, 30, 62, 68, 94, 100)\n\n$tmpNames = @()\n$k = 0\nforeach ($off in $boundaries) {\n    $k = $k + 1\n    $tmpName = \"TempSplitBM\" + $k\n    $tmpNames += $tmpName\n    $pos = $start + $off\n    $splitRange = $d.Range($pos, $pos)\n    $d.Bookmarks.Add($tmpName, $splitRange)\n}\nforeach ($nm in $tmpNames) {\n    $d.Bookmarks($nm).Delete()\n}\n\n# Wrap the first literal \"<br />\" (offsets 24-30) in the real bookmark.\n$brStart = $start + 24\n$brEnd = $start + 30\n$bmRange = $d.Range($brStart, $brEnd)\n$d.Bookmarks.Add(\"__DdeLink__256_169685483\", $bmRange)\n"}
